# Update the "Översikt HÖRBY" workbook:
#  1. Bump the "Förändrad" (C) date from 2023-09-19 to 2023-09-20 for all
#     existing data rows (2-121).
#  2. Give row 121 an explicit 15pt row height (it previously had none).
#  3. Append five brand-new rows (122-126) of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Advance column C (Förändrad) by one day for every existing data row
# ---------------------------------------------------------------------
for ($r = 2; $r -le 121; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value() = $cell.Value().AddDays(1)
}

# ---------------------------------------------------------------------
# 2. Row 121 now gets an explicit (custom) row height of 15, matching the
#    rest of the sheet's rows.
# ---------------------------------------------------------------------
$ws.Rows.Item(121).RowHeight = 15

# ---------------------------------------------------------------------
# 3. Append the five new rows reported by the source site for HÖRBY.
# ---------------------------------------------------------------------
$changed = Get-Date -Year 2023 -Month 9 -Day 20 -Hour 0 -Minute 0 -Second 0
$reported = Get-Date -Year 2023 -Month 9 -Day 18 -Hour 0 -Minute 0 -Second 0

$newRows = @(
    @{ Row = 122; A = "A 43879-2023"; F = "Sveaskog"; G = 1.1; Height = $true  },
    @{ Row = 123; A = "A 43932-2023"; F = "";         G = 3;   Height = $true  },
    @{ Row = 124; A = "A 43939-2023"; F = "";         G = 0.6; Height = $true  },
    @{ Row = 125; A = "A 43937-2023"; F = "";         G = 2.6; Height = $true  },
    @{ Row = 126; A = "A 43935-2023"; F = "";         G = 2.9; Height = $false }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value() = $item.A

    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 2).Value() = $reported

    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 3).Value() = $changed

    $ws.Cells.Item($r, 4).Value() = "SKÅNE LÄN"
    $ws.Cells.Item($r, 5).Value() = "HÖRBY"

    if ($item.F -ne "") {
        $ws.Cells.Item($r, 6).Value() = $item.F
    }

    $ws.Cells.Item($r, 7).Value() = $item.G

    $zeroRange = $ws.Range($ws.Cells.Item($r, 8), $ws.Cells.Item($r, 17))
    $zeroRange.Value() = 0

    $ws.Cells.Item($r, 18).Value() = ""
    $ws.Cells.Item($r, 18).WrapText = $true

    if ($item.Height) {
        $ws.Rows.Item($r).RowHeight = 15
    }
}
